# Task Enumeration update: "updated deliverables and task enumeration"
#
# Logical changes to the "Realization" section of the sheet (originally rows 24-31):
#   - B25 "Eagle library Resarch"  -> renamed to "Prebuilt Eagle library Research"
#   - B26 "Eagle library build"     stays as-is
#   - New deliverable breakdown (Device / Package / Symbol) inserted before "Eagle schematic"
#   - B27 "Eagle schematic" shifts down (now row 30), value/number unchanged
#   - New deliverable breakdown (Place / Wire / Naming / Error Checking) inserted
#     before the old "Eagle board" row
#   - B28 "Eagle board" -> renamed to "Eagle Layout" (now row 35)
#   - "Bring Up"/"Test " rows shift down accordingly (now rows 37/38)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the old row 27 ("Eagle schematic"), pushing it (and
# everything below) down to row 30.
$ws.Range("A27:A29").EntireRow.Insert()

# Insert 4 more blank rows above the old "Eagle board" row (now at row 31),
# pushing it down to row 35.
$ws.Range("A31:A34").EntireRow.Insert()

# Fill in the new / renamed cell values in the same order they were authored,
# so newly created shared strings come out in the expected sequence.
$ws.Range("B35").Value = "Eagle Layout"
$ws.Range("C27").Value = "Device"
$ws.Range("C28").Value = "Package "
$ws.Range("C29").Value = "Symbol"
$ws.Range("B25").Value = "Prebuilt Eagle library Research"
$ws.Range("C31").Value = "Place"
$ws.Range("C32").Value = "Wire"
$ws.Range("C33").Value = "Naming"
$ws.Range("C34").Value = "Error Checking"

# Restore the view state: zoomed in on the newly added rows.
$excel.ActiveWindow.Zoom = 190
$ws.Range("D30:F35").Select() | Out-Null
